$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Capture the existing data rows (2..26) using .Formula (avoids a quirk
#    where reading .Value on a blank cell returns a bogus reflection string).
# ---------------------------------------------------------------------------
$rowsData = @()
for ($r = 2; $r -le 26; $r++) {
    $a = $ws.Cells.Item($r, 1).Formula
    $b = $ws.Cells.Item($r, 2).Formula
    $c = $ws.Cells.Item($r, 3).Formula
    $d = $ws.Cells.Item($r, 4).Formula
    $e = $ws.Cells.Item($r, 5).Formula
    $f = $ws.Cells.Item($r, 6).Formula
    $rowsData += ,@($a,$b,$c,$d,$e,$f)
}

# ---------------------------------------------------------------------------
# 2. Clear all existing hyperlinks up front - they do not travel with the
#    cells when the grid gets rewritten below, so everything is rebuilt from
#    scratch at the end once the data is in its final position.
# ---------------------------------------------------------------------------
$ws.Range("F2").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 3. Extend the table by one row: give row 27 the same formatting as the
#    (current, about to move) last data row, then push every existing data
#    row down by one, from the bottom up so nothing gets clobbered.
# ---------------------------------------------------------------------------
$ws.Range("A26:F26").Copy()
$ws.Range("A27:F27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($i = $rowsData.Count - 1; $i -ge 0; $i--) {
    $destRow = $i + 3   # old row ($i+2) now lives at ($i+3)
    $vals = $rowsData[$i]
    $ws.Cells.Item($destRow, 1).Formula = $vals[0]
    $ws.Cells.Item($destRow, 2).Formula = $vals[1]
    $ws.Cells.Item($destRow, 3).Formula = $vals[2]
    $ws.Cells.Item($destRow, 4).Formula = $vals[3]
    # Leading apostrophe forces literal text, so "dd.mm.yyyy" strings whose
    # day part is <= 12 are not silently reinterpreted as a real date.
    $ws.Cells.Item($destRow, 5).Formula = "'" + $vals[4]
    $ws.Cells.Item($destRow, 6).Formula = $vals[5]
}

# ---------------------------------------------------------------------------
# 4. Write the brand-new latest-price row into row 2, reusing row 3's style.
# ---------------------------------------------------------------------------
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(2, 1).Formula = "26"
$ws.Cells.Item(2, 2).Formula = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Cells.Item(2, 3).Formula = "P1020"
$ws.Cells.Item(2, 4).Formula = "265"
$ws.Cells.Item(2, 5).Formula = "'23.08.2025"
$ws.Cells.Item(2, 6).Formula = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-august-2025.pdf"

# ---------------------------------------------------------------------------
# 5. Rebuild all hyperlinks for column F, rows 2..7.
# ---------------------------------------------------------------------------
$links = @{
    2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-august-2025.pdf"
    3 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf"
    4 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-19-august-2025.pdf"
    5 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-14-august-2025.pdf"
    6 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf"
    7 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf"
}
foreach ($r in 2..7) {
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $links[$r])
}

# Adding hyperlinks applies Excel's built-in "Hyperlink" style (underline /
# theme color) to the target cells - restore the plain table styling that
# every other cell in column F already uses.
$ws.Range("A2:A7").Copy()
$ws.Range("F2:F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The leading apostrophe used above to protect the date strings also marks
# the cells with a "quote prefix" style flag - strip it back out again so
# column E ends up using the same plain style as the rest of the table.
$ws.Range("A2:A27").Copy()
$ws.Range("E2:E27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A1").Select()

Write-Host "edit complete"
